# "quick analysis test without electrolysis"
#
# - B28 ("iteration1" label text) -> "iterationHISTORICAL"
# - B21: drop the "=B8" formula, keep the cached value (3) as a plain
#   number, and flag it with the yellow "manual override" fill used
#   elsewhere in the sheet (e.g. B18) so it stands out from formula cells
# - B33 ("Limit investment to potentials"): FALSE -> TRUE
# - Move the on-screen selection from B33 to B22 (where the analyst was
#   last looking while doing this quick test)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")
$ws.Activate() | Out-Null

# Replace the formula in B21 with its static value and highlight it.
$ws.Range("B21").Value = 3
$ws.Range("B21").Interior.Color = 65535

# Relabel the iteration setting.
$ws.Range("B28").Value = "iterationHISTORICAL"

# Limit investment to potentials -> TRUE (quick test without electrolysis).
$ws.Range("B33").Value = $true

# Update the view's active cell / scroll position.
$ws.Range("B22").Select() | Out-Null
